# Applies the commit "Upload the new file":
#  - Slide 1 title text: "ChatGPT in Medicine" -> "ChatGPT in Big Data"
#  - 4 new "Title Only" slides (2-5) added after slide 1, each with a
#    title placeholder and a free-floating bullet textbox.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1: retitle
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "ChatGPT in Big Data"

# NOTE: this runtime's function parameter binding does not reliably
# bind named (-Param value) arguments, so call positionally:
#   Add-BulletSlide <Index> <Title> <Bullets array>
function Add-BulletSlide {
    param($Index, $Title, $Bullets)

    # layout 11 = ppLayoutTitleOnly ("Title Only") - title placeholder, no body
    $slide = $p.Slides.Add($Index, 11)
    $slide.Shapes.Item(1).TextFrame.TextRange.Text = $Title

    $tb = $slide.Shapes.AddTextbox(1, 72, 144, 432, 144)
    $tb.Name = "TextBox 2"
    $tb.Fill.Visible = $false
    $tb.TextFrame.WordWrap = -1
    $tb.TextFrame.AutoSize = 1

    $tr = $tb.TextFrame.TextRange
    $tr.Text = [string]::Join("`r", $Bullets)
    # leading blank paragraph, like the authored deck
    $tr.InsertBefore("`r")

    for ($i = 2; $i -le ($Bullets.Count + 1); $i++) {
        $para = $tr.Paragraphs($i, 1)
        $para.Font.Size = 16
        $para.Font.Bold = $false
        $para.Font.Name = "Calibri"
    }

    # AutoSize recalculated the height while typing; restore the
    # authored box height (144pt == 1828800 EMU).
    $tb.Height = 144

    return $slide
}

# ---------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------
Add-BulletSlide 2 "ChatGPT and its application in Big Data" @(
    "ChatGPT as a conversational AI tool in Big Data",
    "Improving data analysis and decision-making through ChatGPT",
    "Enhancing customer experiences with ChatGPT in Big Data"
) | Out-Null

# ---------------------------------------------------------------------
# Slide 3
# ---------------------------------------------------------------------
Add-BulletSlide 3 "How ChatGPT aids in Big Data Collection and Analysis" @(
    "ChatGPT for data collection and extraction from unstructured sources",
    "ChatGPT's ability to perform automated data analysis and summarization",
    "ChatGPT's role in data cleansing and pre-processing for effective analysis"
) | Out-Null

# ---------------------------------------------------------------------
# Slide 4
# ---------------------------------------------------------------------
Add-BulletSlide 4 "Benefits of ChatGPT in Big Data Processing" @(
    "ChatGPT's ability to process large amounts of data in real-time",
    "Enhanced accuracy and efficiency through ChatGPT's natural language processing capabilities",
    "ChatGPT's automated data processing and analysis leading to reduced time and operational costs"
) | Out-Null

# ---------------------------------------------------------------------
# Slide 5
# ---------------------------------------------------------------------
Add-BulletSlide 5 "ChatGPT and its Role in Predictive Analytics" @(
    "ChatGPT's ability to identify patterns and trends in Big Data",
    "Utilizing ChatGPT for accurate predictions and forecasting",
    "Improved business outcomes through the integration of ChatGPT in Big Data predictive analytics"
) | Out-Null

Write-Output "Slides: $($p.Slides.Count)"
